{"js": "// Old -> new text pairs, listed in the same order they appear in the\n// document body (the date line, then each \"NNN\u00d7N=\" table-cell problem).\nconst pairs = [\n  [\"2025-06-15 Sunday\", \"2025-06-16 Monday\"],\n  [\"142\u00d73=\", \"341\u00d78=\"],\n  [\"866\u00d79=\", \"688\u00d76=\"],\n  [\"603\u00d75=\", \"179\u00d72=\"],\n  [\"698\u00d72=\", \"639\u00d78=\"],\n  [\"881\u00d72=\", \"833\u00d72=\"],\n  [\"133\u00d75=\", \"412\u00d76=\"],\n  [\"129\u00d72=\", \"743\u00d75=\"],\n  [\"146\u00d77=\", \"983\u00d72=\"],\n  [\"722\u00d72=\", \"236\u00d72=\"],\n  [\"587\u00d76=\", \"527\u00d76=\"],\n  [\"953\u00d74=\", \"961\u00d76=\"],\n  [\"155\u00d79=\", \"657\u00d72=\"],\n  [\"987\u00d79=\", \"733\u00d77=\"],\n  [\"954\u00d73=\", \"552\u00d72=\"],\n  [\"145\u00d72=\", \"884\u00d75=\"],\n  [\"708\u00d79=\", \"372\u00d78=\"],\n  [\"624\u00d74=\", \"543\u00d79=\"],\n  [\"350\u00d72=\", \"607\u00d74=\"],\n  [\"388\u00d75=\", \"334\u00d76=\"],\n  [\"778\u00d74=\", \"846\u00d74=\"],\n  [\"986\u00d73=\", \"896\u00d73=\"],\n  [\"550\u00d74=\", \"627\u00d72=\"],\n  [\"465\u00d76=\", \"987\u00d79=\"],\n  [\"401\u00d72=\", \"914\u00d73=\"],\n  [\"908\u00d76=\", \"892\u00d78=\"],\n];\n\n// Each pair is searched for and replaced one at a time (search -> sync ->\n// replace) rather than batching every search up front: one of the new\n// values (\"987\u00d79=\") is identical to an old value used earlier in the list,\n// so resolving + replacing sequentially guarantees every search sees the\n// document state as of that step instead of a stale snapshot.\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-06-15 Sunday\", \"2025-06-16 Monday\"),\n    @(\"142\u00d73=\", \"341\u00d78=\"),\n    @(\"866\u00d79=\", \"688\u00d76=\"),\n    @(\"603\u00d75=\", \"179\u00d72=\"),\n    @(\"698\u00d72=\", \"639\u00d78=\"),\n    @(\"881\u00d72=\", \"833\u00d72=\"),\n    @(\"133\u00d75=\", \"412\u00d76=\"),\n    @(\"129\u00d72=\", \"743\u00d75=\"),\n    @(\"146\u00d77=\", \"983\u00d72=\"),\n    @(\"722\u00d72=\", \"236\u00d72=\"),\n    @(\"587\u00d76=\", \"527\u00d76=\"),\n    @(\"953\u00d74=\", \"961\u00d76=\"),\n    @(\"155\u00d79=\", \"657\u00d72=\"),\n    @(\"987\u00d79=\", \"733\u00d77=\"),\n    @(\"954\u00d73=\", \"552\u00d72=\"),\n    @(\"145\u00d72=\", \"884\u00d75=\"),\n    @(\"708\u00d79=\", \"372\u00d78=\"),\n    @(\"624\u00d74=\", \"543\u00d79=\"),\n    @(\"350\u00d72=\", \"607\u00d74=\"),\n    @(\"388\u00d75=\", \"334\u00d76=\"),\n    @(\"778\u00d74=\", \"846\u00d74=\"),\n    @(\"986\u00d73=\", \"896\u00d73=\"),\n    @(\"550\u00d74=\", \"627\u00d72=\"),\n    @(\"465\u00d76=\", \"987\u00d79=\"),\n    @(\"401\u00d72=\", \"914\u00d73=\"),\n    @(\"908\u00d76=\", \"892\u00d78=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n\nWrite-Output \"done\"\n"}
